$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 86560965.73999999
$ws.Range("P2").Value = 773816127.41
$ws.Range("Q2").Value = 673010630.39
$ws.Range("R2").Value = 29.138228815
$ws.Range("S2").Value = 555613136.96
$ws.Range("T2").Value = 555613136.96
$ws.Range("U2").Value = 29.6781601625
$ws.Range("V2").Value = 24499989.81
$ws.Range("W2").Value = 48144239.46
$ws.Range("X2").Value = 6561897.44
$ws.Range("Y2").Value = 95820603.75
$ws.Range("Z2").Value = 98378140.88
$ws.Range("AA2").Value = 13974418.19
$ws.Range("AG2").Value = 7321230.35
$ws.Range("AP2").Value = 29.6684964475
$ws.Range("AQ2").Value = 46.588293496325
$ws.Range("AR2").Value = 51.613156894881
$ws.Range("AS2").Value = 84117681.34
$ws.Range("AT2").Value = 42.167805114972
